$d = $word.ActiveDocument

# ============================================================
# Edit 1: "...valid stock dominated in USD..." ->
#         "...valid stock denominated in USD..."
# The original run gets split into three runs so that the
# corrected word "denominated" sits in its own run.
# ============================================================

# Locate the word "dominated" precisely (whole word only, so we
# don't accidentally match a substring of "denominated" etc.).
$findWord = $d.Content
$found = $findWord.Find.Execute("dominated", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'dominated'"
}
$wordStart = $findWord.Start
$wordEnd = $findWord.End

# Replace just that word's text in place.
$rWord = $d.Range($wordStart, $wordEnd)
$rWord.Text = "denominated"

# The replacement word is 2 characters longer than the original,
# so its own range is now ($wordStart, $wordStart + 11).
$newWordEnd = $wordStart + 11
$rDenominated = $d.Range($wordStart, $newWordEnd)

# Toggling a character formatting property and reverting it forces
# the run to be split out from its neighbours without altering the
# visible formatting.
$rDenominated.Font.Bold = $true
$rDenominated.Font.Bold = $false

# Re-locate "must" (now shifted by +2 characters) so the pre-existing
# run boundary between "must" and " contain only..." is preserved
# rather than being swallowed into the newly edited text.
$tailRange = $d.Range($newWordEnd, $d.Content.End)
$foundMust = $tailRange.Find.Execute("must", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMust) {
    $rMust = $d.Range($tailRange.Start, $tailRange.End)
    $rMust.Font.Bold = $true
    $rMust.Font.Bold = $false
}

# ============================================================
# Edit 2: merge the separate "." and "  " runs (after "directly
# for discussion") into a single run containing ".  ".
# ============================================================

$findPeriod = $d.Content
$foundPeriod = $findPeriod.Find.Execute("discussion.  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPeriod) {
    throw "Could not find 'discussion.  '"
}
$spacesStart = $findPeriod.End - 2
$spacesEnd = $findPeriod.End

# Run a (no-visible-change) Find/Replace scoped to exactly the
# trailing two-space run; this forces the engine to re-merge it
# with the immediately preceding "." run (same formatting) without
# touching the unrelated " directly for discussion" run before it.
$rSpaces = $d.Range($spacesStart, $spacesEnd)
$rSpaces.Find.Execute("  ", $true, $false, $false, $false, $false, $true, 1, $false, "  ", 2) | Out-Null
